$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff / Handback datetimes for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 08:43:48"
$wsZhCn.Range("H2").Value = "2016-03-11 08:44:05"

# de-de sheet: update Correspond Handoff / Handback datetimes for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 08:43:51"
$wsDeDe.Range("H2").Value = "2016-03-11 08:44:11"
